$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1
    3 = 4
    4 = 9
    5 = 4
    6 = 7
    7 = 6
    8 = 9
    9 = 7
    10 = 9
    11 = 2
    12 = 9
    13 = 5
    14 = 5
    15 = 3
    16 = 3
    17 = 2
    18 = 8
    19 = 1
    20 = 1
    21 = 3
    22 = 5
    23 = 4
    24 = 3
    25 = 3
    26 = 0
    27 = 4
    28 = 7
    29 = 2
    30 = 9
    31 = 3
    32 = 2
    33 = 2
    34 = 5
    35 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
